# Auto-generated edit script for 'Add data for 2025-08-06'
# Updates 2025 (column L) cumulative counts, plus a few retroactive
# corrections to columns I/K (2022/2024) across Citywide Totals,
# By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 4022
$ws.Range("L3").Value = 4236
$ws.Range("I4").Value = 1841
$ws.Range("K4").Value = 1779
$ws.Range("L4").Value = 1048
$ws.Range("L5").Value = 244
$ws.Range("L6").Value = 3649
$ws.Range("I7").Value = 26310
$ws.Range("K7").Value = 27571
$ws.Range("L7").Value = 13199

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L2").Value = 46
$ws.Range("L3").Value = 40
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 252
$ws.Range("L3").Value = 290
$ws.Range("L4").Value = 61
$ws.Range("L6").Value = 235
$ws.Range("L7").Value = 866

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 96
$ws.Range("L3").Value = 116
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 296

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 207
$ws.Range("L7").Value = 618

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 57
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 161
$ws.Range("L6").Value = 133
$ws.Range("L7").Value = 484

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 56
$ws.Range("L7").Value = 216

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 110
$ws.Range("L7").Value = 437
$ws.Range("L8").Value = 866
$ws.Range("L10").Value = 86
$ws.Range("L11").Value = 214
$ws.Range("L14").Value = 64
$ws.Range("L15").Value = 97
$ws.Range("L19").Value = 377
$ws.Range("L20").Value = 330
$ws.Range("L23").Value = 141
$ws.Range("L29").Value = 742
$ws.Range("L31").Value = 126
$ws.Range("L33").Value = 618
$ws.Range("L37").Value = 484
$ws.Range("L40").Value = 37
$ws.Range("I41").Value = 118
$ws.Range("L42").Value = 422
$ws.Range("L44").Value = 94
$ws.Range("L50").Value = 65
$ws.Range("L52").Value = 265
$ws.Range("L53").Value = 150
$ws.Range("L54").Value = 274
$ws.Range("L55").Value = 126
$ws.Range("K63").Value = 168
$ws.Range("L63").Value = 42
$ws.Range("L67").Value = 463
$ws.Range("L74").Value = 12
$ws.Range("L78").Value = 167
$ws.Range("L79").Value = 347
$ws.Range("L83").Value = 296
$ws.Range("L84").Value = 128
$ws.Range("L85").Value = 689
$ws.Range("L94").Value = 166
$ws.Range("L95").Value = 181
$ws.Range("L97").Value = 110
$ws.Range("L98").Value = 74
$ws.Range("L99").Value = 216
$ws.Range("I101").Value = 26310
$ws.Range("K101").Value = 27571
$ws.Range("L101").Value = 13199

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 131
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 463

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 128

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 274

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 227
$ws.Range("L3").Value = 280
$ws.Range("L5").Value = 13
$ws.Range("L6").Value = 188
$ws.Range("L7").Value = 742

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 132
$ws.Range("L7").Value = 377

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 125
$ws.Range("L7").Value = 422

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 51
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L2").Value = 42
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 37
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 141

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 116
$ws.Range("L7").Value = 347

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L4").Value = 32
$ws.Range("L6").Value = 87
$ws.Range("L7").Value = 330

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 142
$ws.Range("L4").Value = 31
$ws.Range("L6").Value = 119
$ws.Range("L7").Value = 437

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 37
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 166

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 36
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L6").Value = 52
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 110

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 212
$ws.Range("L3").Value = 275
$ws.Range("L7").Value = 689

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 83
$ws.Range("L6").Value = 71
$ws.Range("L7").Value = 265

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 12
